$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.878.93"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.59"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4771"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3929"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.85"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07986"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.010"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.74"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.884.72"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.018"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.173"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06695"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001045"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.874.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.490"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.67%  "

$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.330"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.097.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.75"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.095"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.456"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9722"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09495"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.628"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.324"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.348"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06060"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.72%  "

$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.203"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.60%  "

$ws.Range("E40").Value = "  -1.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.008"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5950"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1889"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.85%  "

$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.255"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5662"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.16"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.920"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.314"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06780"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.12"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.01%  "
